$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pin description: ENOP -> ENOP TIM3_CH2
$ws.Range("D5").Value = "ENOP TIM3_CH2"

# Match the recorded active cell selection after the edit
$ws.Range("D16").Select() | Out-Null
